$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.830.46'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = '2.600.60'
$ws.Range("E3").Value = '  -1.44%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.59'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.00%  '

$ws.Range("D9").Value = '2.602.54'
$ws.Range("E9").Value = '  -1.36%  '

$ws.Range("E10").Value = '  -2.67%  '

$ws.Range("E11").Value = '  +0.42%  '

$ws.Range("E12").Value = '  +0.59%  '

$ws.Range("E13").Value = '  -1.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.06'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.80%  '

$ws.Range("D15").Value = '3.077.87'

$ws.Range("E16").Value = '  -2.12%  '

$ws.Range("D17").Value = '66.950.40'
$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("D18").Value = '2.605.09'
$ws.Range("E18").Value = '  -2.34%  '

$ws.Range("E19").Value = '  -3.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.77'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.28'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.21%  '

$ws.Range("E22").Value = '  -2.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.60'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.42'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -5.10%  '

$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("E26").Value = '  -4.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '69.05'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.02%  '

$ws.Range("D28").Value = '2.740.64'
$ws.Range("E28").Value = '  -1.38%  '

$ws.Range("E29").Value = '  -0.18%  '

$ws.Range("D30").Value = '0.0₃0984'
$ws.Range("E30").Value = '  -3.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '538.11'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.10'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.93%  '

$ws.Range("E33").Value = '  -3.68%  '

$ws.Range("E34").Value = '  -3.15%  '

$ws.Range("E35").Value = '  -0.27%  '

$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("E37").Value = '  -4.44%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.91'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.80'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.16%  '

$ws.Range("E40").Value = '  -2.32%  '

$ws.Range("E41").Value = '  +1.96%  '

$ws.Range("E42").Value = '  -0.65%  '

$ws.Range("E43").Value = '  -3.14%  '

$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("E45").Value = '  -4.74%  '

$ws.Range("D46").Value = '0.0₆0292'
$ws.Range("E46").Value = '  -1.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.01'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.00%  '

$ws.Range("E48").Value = '  -3.41%  '

$ws.Range("E49").Value = '  -2.75%  '

$ws.Range("E50").Value = '  -1.73%  '

$ws.Range("E51").Value = '  -1.29%  '
